$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Regenerate merged AHB headers: "_old" columns become "_FV2410",
#    "_new" columns become "_FV2504". The "diff" column (K) is untouched.
$oldHeaders = @("Segmentname_old","Segmentgruppe_old","Segment_old","Datenelement_old","Segment ID_old","Code_old","Qualifier_old","Beschreibung_old","Bedingungsausdruck_old","Bedingung_old")
$newHeaders = @("Segmentname_new","Segmentgruppe_new","Segment_new","Datenelement_new","Segment ID_new","Code_new","Qualifier_new","Beschreibung_new","Bedingungsausdruck_new","Bedingung_new")

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = ($oldHeaders[$i] -replace '_old$', '_FV2410')
}
for ($i = 0; $i -lt $newHeaders.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = ($newHeaders[$i] -replace '_new$', '_FV2504')
}

# 2) Turn the used range into a real table (Table1), keeping the header
#    row's existing formatting (bold / fill / border) intact. A scratch
#    row is used to stash+restore that formatting across the ClearFormats
#    call that Excel needs before laying the table header styling down,
#    so neither styles.xml's cellXfs/dxfs nor the header look actually
#    change.
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A100:U100")

$headerRange.Copy()
$scratch.PasteSpecial(-4122) | Out-Null

$headerRange.ClearFormats()
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U78"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$scratch.Copy()
$headerRange.PasteSpecial(-4122) | Out-Null
$scratch.ClearFormats()
$scratch.ClearContents()

# 3) Freeze the header row.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Output "Regenerated merged AHB headers, added Table1, froze header row."
